$d = $word.ActiveDocument

$pairs = @(
    @("585÷7=", "614÷5="),
    @("419÷8=", "127÷2="),
    @("997÷7=", "681÷4="),
    @("790÷4=", "130÷6="),
    @("268÷3=", "773÷3="),
    @("126÷5=", "699÷5="),
    @("635÷6=", "669÷4="),
    @("782÷7=", "818÷7="),
    @("431÷9=", "899÷8="),
    @("562÷3=", "288÷9="),
    @("882÷3=", "951÷9="),
    @("984÷8=", "774÷3="),
    @("504÷5=", "461÷9="),
    @("577÷5=", "478÷4="),
    @("675÷7=", "574÷8="),
    @("827÷8=", "935÷9="),
    @("868÷3=", "181÷6="),
    @("745÷2=", "849÷9="),
    @("752÷5=", "169÷8="),
    @("642÷8=", "523÷6="),
    @("744÷6=", "986÷3="),
    @("706÷2=", "707÷8="),
    @("109÷2=", "761÷5="),
    @("119÷7=", "238÷3="),
    @("477÷9=", "224÷2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
